$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (song-index numbers) to fit the new, longer data.
$ws.Columns.Item(1).ColumnWidth = 51.8

# New tracks added on 2018-11-18, appended after the existing 46 rows
# (excel rows 2..47 already hold index 0..45, so the new ones continue
# the A-column numbering at 46..50 and go into rows 48..52).
$newSongs = @(
    @{ Index = 46; Title = "It Ain_t Me";                                                  Author = "Kygo, Selena Gomez";          Date = "18-11-2018" },
    @{ Index = 47; Title = "Don_t Let Me Down";                                             Author = "The Chainsmokers, Da";        Date = "18-11-2018" },
    @{ Index = 48; Title = "Complicated_Dimitri_Vegas_Like_Mike_D_[500kbps_M4A]";            Author = "Unknown";                     Date = "18-11-2018" },
    @{ Index = 49; Title = "2U";                                                            Author = "David Guetta, Justin Bieber"; Date = "18-11-2018" },
    @{ Index = 50; Title = "Attention";                                                     Author = "Charlie Puth";                Date = "18-11-2018" }
)

$firstNewRow = 48
$lastNewRow = $firstNewRow + $newSongs.Count - 1

# Carry the styled (bold/border/centered) look of column A down into the
# new rows by copying the format from the last existing data row.
$ws.Range("A47").Copy() | Out-Null
$ws.Range("A$firstNewRow`:A$lastNewRow").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $newSongs.Count; $i++) {
    $row = $firstNewRow + $i
    $song = $newSongs[$i]
    $ws.Cells.Item($row, 1).Value = $song.Index
    $ws.Cells.Item($row, 2).Value = $song.Title
    $ws.Cells.Item($row, 3).Value = $song.Author
    $ws.Cells.Item($row, 4).Value = $song.Date
}
